# Kumiko Example.xlsx — FileWriter / GUI save-path update
#
# Header row: "YAAAAAA" -> "Project Name!" and "Pieces" -> "Piece #".
# The "Piece #" column (B) is now filled with a running piece count
# (1..10), replacing the old scratch values that used to live in the
# "Material" column (D2:D3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Header relabel -------------------------------------------------
$ws.Range("A1").Value = "Project Name!"
$ws.Range("B1").Value = "Piece #"

# -- Drop the old leftover Material test values ----------------------
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

# -- Populate the Piece # column with rows 2-11 -----------------------
$pieceNumbers = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)
for ($i = 0; $i -lt $pieceNumbers.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $pieceNumbers[$i]
}

# -- Row heights: first block (header + initial pieces) shrinks a bit,
#    the newly-appended rows keep the sheet's normal default height.
$ws.Rows("1:7").RowHeight = 13.2
$ws.Rows("8:11").RowHeight = 15.75

# -- Widen column A (Project Name!) so the longer header fits --------
$ws.Columns("A:A").ColumnWidth = 12.8

# -- Leave the selection where the author left it ---------------------
$ws.Range("D5").Select() | Out-Null
